$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H5").Value = '[''shore_mon_hr'', ''shore_mon_fut'', ''shore_mon'']'
$ws.Range("E8").Value = 'string'
$ws.Range("H8").Value = '[''sed_class'', ''shore_mon'', ''world_pop'', ''shore_mon_drivers'', ''shore_mon_fut'', ''world_gdp'']'
$ws.Range("E9").Value = 'string'
$ws.Range("H9").Value = '[''sed_class'', ''shore_mon'', ''world_pop'', ''shore_mon_drivers'', ''shore_mon_fut'', ''world_gdp'']'
$ws.Range("E10").Value = 'string'
$ws.Range("H10").Value = '[''shore_mon_fut'', ''shore_mon_drivers'', ''shore_mon'']'
$ws.Range("E13").Value = 'string'
$ws.Range("H13").Value = '[''shore_mon_hr'', ''shore_mon_drivers'']'
$ws.Range("H17").Value = '[''sed_class'', ''shore_mon_hr'', ''shore_mon'', ''world_pop'', ''shore_mon_drivers'', ''shore_mon_fut'', ''world_gdp'']'
$ws.Range("H19").Value = '[''sed_class'', ''shore_mon_hr'', ''shore_mon'', ''world_pop'', ''shore_mon_drivers'', ''shore_mon_fut'', ''world_gdp'']'
$ws.Range("E25").Value = 'string'
$ws.Range("E27").Value = 'string'
$ws.Range("H27").Value = '[''sed_class'', ''shore_mon_hr'', ''shore_mon'', ''world_pop'', ''shore_mon_drivers'', ''shore_mon_fut'', ''world_gdp'']'
$ws.Range("I27").Value = '[b''LINESTRING (-72.98252423386046 -54.45955389684386, -72.99521062694167 -54.448286932134025)'', b''LINESTRING (-74.3863095545 -50.3776589451, -74.39562299239999 -50.3874558377)'', b''LINESTRING (-71.624166 10.990024, -71.60742 10.996413)'', b''LINESTRING (47.885022 29.319149, 47.879413 29.33403)'', b''LINESTRING (-74.382468591 -50.3791437735, -74.3917820288 -50.3889403594)'']'
$ws.Range("E28").Value = 'string'
$ws.Range("H28").Value = '[''sed_class'', ''shore_mon_hr'', ''shore_mon'', ''world_pop'', ''shore_mon_drivers'', ''shore_mon_fut'', ''world_gdp'']'
$ws.Range("I28").Value = '[b''BOX_028_183_1'', b''BOX_028_183_0'', b''BOX_028_000_0'', b''BOX_117_067_110'', b''BOX_145_168_39'']'
$ws.Range("H29").Value = '[''shore_mon_fut'', ''shore_mon_drivers'', ''shore_mon'']'
$ws.Range("E31").Value = 'string'
$ws.Range("E35").Value = 'string'
$ws.Range("E36").Value = 'string'
$ws.Range("E37").Value = 'string'
$ws.Range("H38").Value = '[''shore_mon_hr'', ''shore_mon'']'
$ws.Range("E43").Value = 'string'
$ws.Range("E46").Value = 'string'
$ws.Range("H48").Value = '[''shore_mon_hr'', ''shore_mon'']'
